$d = $word.ActiveDocument

# --- Change 1: "Material necessário" list item - merge spell-checked runs ---
# Before: proofErr-wrapped "Microcontrolador" + " " + "Steps" + ";" (4 runs, 3 proofErr markers)
# After: single run "Microcontrolador Steps;"
$p12 = $d.Paragraphs.Item(12)
$r12 = $p12.Range
$r12.Delete()
$r12.InsertBefore("Microcontrolador Steps;`r")

# --- Change 2: restructure the "Aplicação"/bookmark/"Exemplos:"/URL block ---
# 2a: "Aplicação" -> "Outras abordagens"
$d.Paragraphs.Item(19).Range.Text = "Outras abordagens"

# 2b: remove the bookmark-only paragraph (merges away, leaving two empty paragraphs)
$d.Paragraphs.Item(20).Range.Delete()

# 2c: duplicate the now-empty paragraph so there are two blank paragraphs
$pBlank = $d.Paragraphs.Item(20)
$rBlank = $pBlank.Range.Duplicate
$rBlank.Collapse(0)
$rBlank.InsertAfter("`r")

# 2d: remove old "Exemplos:" + URL paragraphs, replace with "Veja mais:" + URL
$pLastBlank = $d.Paragraphs.Item(21)
$pUrl = $d.Paragraphs.Item(23)
$full = $d.Range($pLastBlank.Range.End, $pUrl.Range.End)
$full.Delete()

$insertPoint = $d.Paragraphs.Item(21).Range.Duplicate
$insertPoint.Collapse(0)
$insertPoint.InsertAfter("`rVeja mais:`rhttp://www.matematicadidatica.com.br/GeometriaCalculoAreaFigurasPlanas.aspx")

# 2e: move the _GoBack bookmark in between "Veja mais" and ":"
$pVeja = $d.Paragraphs.Item(22)
$bmPos = $pVeja.Range.Duplicate
$bmPos.Collapse(1)
$bmPos.MoveEnd(1, 9) | Out-Null
$bmPos.Collapse(0)
$d.Bookmarks.Add("_GoBack", $bmPos)
